# Split the single "Terms Typically Offered" column (D) into four columns:
# D=Corequisites, E=Concurrent, F=Recommended, G=Terms Typically Offered.
# This inserts three new blank columns before the existing column D, which
# pushes the current "Terms Typically Offered" data from D into G automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:F").Insert()

# New header row.
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"

# Default every data row's new Corequisites / Concurrent / Recommended cells
# to "NA"; specific rows are corrected afterwards below.
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 4).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
}

# Row 11 (ECON 313): the old "Terms Typically Offered" cell had a corequisite
# note glued onto the end ("F, W, SPCorequisite: ECON 311."). Split it: the
# corequisite goes to the Corequisites column, and the terms-offered text
# (now in column G after the insert) is corrected to just the terms.
$ws.Cells.Item(11, 4).Value = "ECON 311."
$ws.Cells.Item(11, 7).Value = "F, W, SP"

# Row 31 (ECON 434): the old Prerequisites cell had a "Recommended: ECON 312."
# clause appended. Move that into the new Recommended column and trim the
# Prerequisites text back down.
$ws.Cells.Item(31, 3).Value = "ECON 311."
$ws.Cells.Item(31, 6).Value = "ECON 312."
$ws.Cells.Item(31, 7).Value = "SP "

# Row 37 (ECON 464): the old Prerequisites cell had a "Corequisite: ECON 460."
# clause appended. Move that into the new Corequisites column and trim the
# Prerequisites text back down.
$ws.Cells.Item(37, 3).Value = "Senior standing and two 400-level ECON courses other than ECON 460."
$ws.Cells.Item(37, 4).Value = "ECON 460."
$ws.Cells.Item(37, 7).Value = "F, SP "
